# Splitting ratio is now -3dB
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the LOG10-based splitting-ratio formulas in D18/E18 with a
# static -3 dB value each (previously =LOG10(2/3)*10 and =LOG10(1/3)*10).
$ws.Range("D18").Value = -3
$ws.Range("E18").Value = -3

# Move/restore the saved selection to E19, matching the workbook's view state.
$ws.Range("E19").Select()
